$d = $word.ActiveDocument
$newText = "Du deltar i en världsomspännande kampanj för att observera och rapportera de svagaste synliga stjärnorna, som ett mått på ljusföroreningarna på orten. Genom att hitta och observera Perseus konstellation på natthimlen kan folk i hela världen lära sig hur belysningen i våra samhällen och omgivningar bidrar till ljusföroreningar. Era bidrag till online-databasen hjälper till att dokumentera den synliga natthimlens över hela världen."

# The four "Kampanjdatum for Perseus ..." banner paragraphs and the one
# "Du deltar i en varldsomspannande..." intro paragraph are all being
# collapsed into a single plain run carrying the (new) intro copy.
# Identify them by their distinctive old text (rather than a fixed
# paragraph index) so the script is resilient to renumbering; walk
# back-to-front so edits never shift the index of a not-yet-visited
# paragraph.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Kampanjdatum f*" -or $t -like "Du deltar i en v*") {
        $r = $p.Range
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.Delete()
        $r2.InsertAfter($newText)
    }
}
